# Update NATMI LR-pair output (Sema4c-Plxnb2) with refreshed TPM-derived values.
# The underlying ligand (Sema4c) and receptor (Plxnb2) average/total expression
# values for the "ECs" cluster were recomputed from new TPM input; every
# dependent specificity / edge-weight column is refreshed to stay consistent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.06647866666667
$ws.Range("H2").Value = 33.199436
$ws.Range("I2").Value = 0.4391813340952289
$ws.Range("J2").Value = 0.4391813340952289
$ws.Range("M2").Value = 2.325008666666667
$ws.Range("N2").Value = 6.975026
$ws.Range("O2").Value = 0.05445297772988467
$ws.Range("P2").Value = 0.05445297772988466
$ws.Range("Q2").Value = 25.72965880948178
$ws.Range("R2").Value = 231.566929285336
$ws.Range("S2").Value = 0.02391473140486854
$ws.Range("T2").Value = 0.02391473140486853
$ws.Range("G3").Value = 11.06647866666667
$ws.Range("H3").Value = 33.199436
$ws.Range("I3").Value = 0.4391813340952289
$ws.Range("J3").Value = 0.4391813340952289
$ws.Range("O3").Value = 0.4529132218878514
$ws.Range("P3").Value = 0.4529132218878514
$ws.Range("Q3").Value = 214.0067110247676
$ws.Range("R3").Value = 1926.060399222908
$ws.Range("S3").Value = 0.198911033018075
$ws.Range("T3").Value = 0.198911033018075
$ws.Range("G4").Value = 11.06647866666667
$ws.Range("H4").Value = 33.199436
$ws.Range("I4").Value = 0.4391813340952289
$ws.Range("J4").Value = 0.4391813340952289
$ws.Range("O4").Value = 0.492633800382264
$ws.Range("P4").Value = 0.492633800382264
$ws.Range("Q4").Value = 232.7751416043791
$ws.Range("R4").Value = 2094.976274439412
$ws.Range("S4").Value = 0.2163555696722854
$ws.Range("T4").Value = 0.2163555696722853
$ws.Range("I5").Value = 0.4357622493260503
$ws.Range("J5").Value = 0.4357622493260503
$ws.Range("M5").Value = 2.325008666666667
$ws.Range("N5").Value = 6.975026
$ws.Range("O5").Value = 0.05445297772988467
$ws.Range("P5").Value = 0.05445297772988466
$ws.Range("Q5").Value = 25.52935001281378
$ws.Range("R5").Value = 229.764150115324
$ws.Range("S5").Value = 0.02372855205807587
$ws.Range("T5").Value = 0.02372855205807586
$ws.Range("I6").Value = 0.4357622493260503
$ws.Range("J6").Value = 0.4357622493260503
$ws.Range("O6").Value = 0.4529132218878514
$ws.Range("P6").Value = 0.4529132218878514
$ws.Range("S6").Value = 0.1973624843193587
$ws.Range("T6").Value = 0.1973624843193586
$ws.Range("I7").Value = 0.4357622493260503
$ws.Range("J7").Value = 0.4357622493260503
$ws.Range("O7").Value = 0.492633800382264
$ws.Range("P7").Value = 0.492633800382264
$ws.Range("S7").Value = 0.2146712129486158
$ws.Range("T7").Value = 0.2146712129486158
$ws.Range("I8").Value = 0.1250564165787209
$ws.Range("J8").Value = 0.1250564165787209
$ws.Range("M8").Value = 2.325008666666667
$ws.Range("N8").Value = 6.975026
$ws.Range("O8").Value = 0.05445297772988467
$ws.Range("P8").Value = 0.05445297772988466
$ws.Range("Q8").Value = 7.326492910122666
$ws.Range("R8").Value = 65.93843619110399
$ws.Range("S8").Value = 0.006809694266940267
$ws.Range("T8").Value = 0.006809694266940266
$ws.Range("I9").Value = 0.1250564165787209
$ws.Range("J9").Value = 0.1250564165787209
$ws.Range("O9").Value = 0.4529132218878514
$ws.Range("P9").Value = 0.4529132218878514
$ws.Range("S9").Value = 0.05663970455041778
$ws.Range("T9").Value = 0.05663970455041777
$ws.Range("I10").Value = 0.1250564165787209
$ws.Range("J10").Value = 0.1250564165787209
$ws.Range("O10").Value = 0.492633800382264
$ws.Range("P10").Value = 0.492633800382264
$ws.Range("S10").Value = 0.06160701776136282
$ws.Range("T10").Value = 0.06160701776136281
